{"js": "// Replace the two-digit-division answer strings in the table cells with\n// their updated values. Each old value is unique within the document, so a\n// simple body-wide search-and-replace for each pair is unambiguous.\nconst replacements = [\n  [\"18\u00f77=2, 4\", \"49\u00f78=6, 1\"],\n  [\"11\u00f77=1, 4\", \"38\u00f77=5, 3\"],\n  [\"26\u00f76=4, 2\", \"82\u00f76=13, 4\"],\n  [\"12\u00f75=2, 2\", \"33\u00f73=11, 0\"],\n  [\"37\u00f73=12, 1\", \"77\u00f78=9, 5\"],\n  [\"89\u00f76=14, 5\", \"36\u00f78=4, 4\"],\n  [\"98\u00f74=24, 2\", \"85\u00f75=17, 0\"],\n  [\"93\u00f77=13, 2\", \"58\u00f79=6, 4\"],\n  [\"31\u00f78=3, 7\", \"43\u00f79=4, 7\"],\n  [\"64\u00f78=8, 0\", \"19\u00f78=2, 3\"],\n  [\"78\u00f74=19, 2\", \"69\u00f74=17, 1\"],\n  [\"54\u00f77=7, 5\", \"52\u00f78=6, 4\"],\n  [\"87\u00f75=17, 2\", \"60\u00f73=20, 0\"],\n  [\"54\u00f79=6, 0\", \"62\u00f75=12, 2\"],\n  [\"15\u00f77=2, 1\", \"70\u00f78=8, 6\"],\n  [\"70\u00f75=14, 0\", \"33\u00f78=4, 1\"],\n  [\"32\u00f78=4, 0\", \"44\u00f72=22, 0\"],\n  [\"15\u00f75=3, 0\", \"52\u00f72=26, 0\"],\n  [\"82\u00f79=9, 1\", \"28\u00f75=5, 3\"],\n  [\"89\u00f78=11, 1\", \"66\u00f75=13, 1\"],\n  [\"55\u00f73=18, 1\", \"65\u00f73=21, 2\"],\n  [\"79\u00f78=9, 7\", \"76\u00f79=8, 4\"],\n  [\"44\u00f78=5, 4\", \"96\u00f72=48, 0\"],\n  [\"48\u00f78=6, 0\", \"75\u00f75=15, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-division answer strings in the table cells with\n# their updated values. Each old value is unique within the document, so a\n# straightforward Find/Replace (wdReplaceAll) for each pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"18\u00f77=2, 4\", \"49\u00f78=6, 1\"),\n    @(\"11\u00f77=1, 4\", \"38\u00f77=5, 3\"),\n    @(\"26\u00f76=4, 2\", \"82\u00f76=13, 4\"),\n    @(\"12\u00f75=2, 2\", \"33\u00f73=11, 0\"),\n    @(\"37\u00f73=12, 1\", \"77\u00f78=9, 5\"),\n    @(\"89\u00f76=14, 5\", \"36\u00f78=4, 4\"),\n    @(\"98\u00f74=24, 2\", \"85\u00f75=17, 0\"),\n    @(\"93\u00f77=13, 2\", \"58\u00f79=6, 4\"),\n    @(\"31\u00f78=3, 7\", \"43\u00f79=4, 7\"),\n    @(\"64\u00f78=8, 0\", \"19\u00f78=2, 3\"),\n    @(\"78\u00f74=19, 2\", \"69\u00f74=17, 1\"),\n    @(\"54\u00f77=7, 5\", \"52\u00f78=6, 4\"),\n    @(\"87\u00f75=17, 2\", \"60\u00f73=20, 0\"),\n    @(\"54\u00f79=6, 0\", \"62\u00f75=12, 2\"),\n    @(\"15\u00f77=2, 1\", \"70\u00f78=8, 6\"),\n    @(\"70\u00f75=14, 0\", \"33\u00f78=4, 1\"),\n    @(\"32\u00f78=4, 0\", \"44\u00f72=22, 0\"),\n    @(\"15\u00f75=3, 0\", \"52\u00f72=26, 0\"),\n    @(\"82\u00f79=9, 1\", \"28\u00f75=5, 3\"),\n    @(\"89\u00f78=11, 1\", \"66\u00f75=13, 1\"),\n    @(\"55\u00f73=18, 1\", \"65\u00f73=21, 2\"),\n    @(\"79\u00f78=9, 7\", \"76\u00f79=8, 4\"),\n    @(\"44\u00f78=5, 4\", \"96\u00f72=48, 0\"),\n    @(\"48\u00f78=6, 0\", \"75\u00f75=15, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards,\n        $false, $false, $find.Forward, $find.Wrap, $false,\n        $find.Replacement.Text, 2\n    )\n}\n"}
